{"js": "// Replace the date line and every \"A\u00d7B=\" multiplication prompt in the\n// table with its new value. Every source string in this document is\n// unique, so a plain case-sensitive search-and-replace (one pass per\n// pair, oldest text -> newest text) is unambiguous and order independent.\nconst replacements = [\n  [\"2025-03-08 Saturday\", \"2025-03-09 Sunday\"],\n  [\"121\u00d76=\", \"117\u00d78=\"],\n  [\"123\u00d76=\", \"853\u00d78=\"],\n  [\"857\u00d79=\", \"527\u00d73=\"],\n  [\"384\u00d75=\", \"972\u00d76=\"],\n  [\"273\u00d75=\", \"135\u00d79=\"],\n  [\"877\u00d75=\", \"347\u00d73=\"],\n  [\"226\u00d75=\", \"465\u00d78=\"],\n  [\"702\u00d74=\", \"510\u00d77=\"],\n  [\"252\u00d72=\", \"405\u00d78=\"],\n  [\"680\u00d73=\", \"916\u00d79=\"],\n  [\"473\u00d78=\", \"107\u00d72=\"],\n  [\"239\u00d79=\", \"977\u00d79=\"],\n  [\"129\u00d72=\", \"830\u00d76=\"],\n  [\"726\u00d79=\", \"584\u00d76=\"],\n  [\"787\u00d72=\", \"951\u00d76=\"],\n  [\"826\u00d77=\", \"354\u00d79=\"],\n  [\"797\u00d76=\", \"935\u00d78=\"],\n  [\"274\u00d72=\", \"401\u00d77=\"],\n  [\"803\u00d75=\", \"667\u00d73=\"],\n  [\"953\u00d74=\", \"654\u00d72=\"],\n  [\"291\u00d73=\", \"473\u00d75=\"],\n  [\"803\u00d73=\", \"373\u00d75=\"],\n  [\"126\u00d79=\", \"215\u00d72=\"],\n  [\"176\u00d79=\", \"487\u00d76=\"],\n  [\"319\u00d78=\", \"359\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A x B=\" multiplication prompt in the\n# table with its new value. Every source string in this document is\n# unique, so a case-sensitive Find/ReplaceAll pass per pair (oldest text\n# -> newest text) is unambiguous and order independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-08 Saturday\", \"2025-03-09 Sunday\"),\n    @(\"121\u00d76=\", \"117\u00d78=\"),\n    @(\"123\u00d76=\", \"853\u00d78=\"),\n    @(\"857\u00d79=\", \"527\u00d73=\"),\n    @(\"384\u00d75=\", \"972\u00d76=\"),\n    @(\"273\u00d75=\", \"135\u00d79=\"),\n    @(\"877\u00d75=\", \"347\u00d73=\"),\n    @(\"226\u00d75=\", \"465\u00d78=\"),\n    @(\"702\u00d74=\", \"510\u00d77=\"),\n    @(\"252\u00d72=\", \"405\u00d78=\"),\n    @(\"680\u00d73=\", \"916\u00d79=\"),\n    @(\"473\u00d78=\", \"107\u00d72=\"),\n    @(\"239\u00d79=\", \"977\u00d79=\"),\n    @(\"129\u00d72=\", \"830\u00d76=\"),\n    @(\"726\u00d79=\", \"584\u00d76=\"),\n    @(\"787\u00d72=\", \"951\u00d76=\"),\n    @(\"826\u00d77=\", \"354\u00d79=\"),\n    @(\"797\u00d76=\", \"935\u00d78=\"),\n    @(\"274\u00d72=\", \"401\u00d77=\"),\n    @(\"803\u00d75=\", \"667\u00d73=\"),\n    @(\"953\u00d74=\", \"654\u00d72=\"),\n    @(\"291\u00d73=\", \"473\u00d75=\"),\n    @(\"803\u00d73=\", \"373\u00d75=\"),\n    @(\"126\u00d79=\", \"215\u00d72=\"),\n    @(\"176\u00d79=\", \"487\u00d76=\"),\n    @(\"319\u00d78=\", \"359\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\n$d.Paragraphs(1).Range.Text\n"}
